$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H18").Value = 1433.3334
$ws_ALC.Range("I18").Value = 1237.625
$ws_ALC.Range("J18").Value = 2999
$ws_ALC.Range("K18").Value = 1237.625
$ws_ALC.Range("L18").Value = 2999
$ws_ALC.Range("M18").Value = -953.625
$ws_ALC.Range("N18").Value = -3567

$ws_ALC.Range("H28").Value = 300.77274
$ws_ALC.Range("I28").Value = 268.85
$ws_ALC.Range("J28").Value = 620
$ws_ALC.Range("K28").Value = 268.85
$ws_ALC.Range("L28").Value = 620
$ws_ALC.Range("M28").Value = 216.15
$ws_ALC.Range("N28").Value = -1590

$ws_ALC.Range("H48").Value = 0
$ws_ALC.Range("I48").Value = 0
$ws_ALC.Range("J48").Value = 0
$ws_ALC.Range("K48").Value = 0
$ws_ALC.Range("L48").ClearContents()
$ws_ALC.Range("N48").Value = 0

$ws_ALC.Range("H56").Value = 0
$ws_ALC.Range("I56").Value = 0
$ws_ALC.Range("J56").Value = 0
$ws_ALC.Range("K56").Value = 0
$ws_ALC.Range("L56").ClearContents()
$ws_ALC.Range("N56").Value = 0

$ws_ALC.Range("H133").Value = 68885
$ws_ALC.Range("I133").Value = 0
$ws_ALC.Range("J133").Value = 68885
$ws_ALC.Range("K133").Value = 0
$ws_ALC.Range("L133").Value = 68885
$ws_ALC.Range("N133").Value = -79005

$ws_ALC.Range("H139").Value = 53467.777
$ws_ALC.Range("I139").Value = 0
$ws_ALC.Range("J139").Value = 53467.777
$ws_ALC.Range("K139").Value = 0
$ws_ALC.Range("L139").Value = 53467.777
$ws_ALC.Range("N139").Value = -63747.777

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H29").Value = 0
$ws_ARM.Range("I29").Value = 0
$ws_ARM.Range("J29").Value = 0
$ws_ARM.Range("K29").Value = 0
$ws_ARM.Range("L29").ClearContents()
$ws_ARM.Range("N29").Value = 0

$ws_ARM.Range("H32").Value = 14142.296
$ws_ARM.Range("I32").Value = 14142.296
$ws_ARM.Range("J32").Value = 0
$ws_ARM.Range("K32").Value = 14142.296
$ws_ARM.Range("L32").Value = 0
$ws_ARM.Range("M32").ClearContents()
$ws_ARM.Range("N32").Value = -13855.296

$ws_ARM.Range("H74").Value = 1293.762
$ws_ARM.Range("I74").Value = 1098.6666
$ws_ARM.Range("J74").Value = 1440.0834
$ws_ARM.Range("K74").Value = 1098.6666
$ws_ARM.Range("L74").Value = 1440.0834
$ws_ARM.Range("M74").Value = -224.6666
$ws_ARM.Range("N74").Value = -3188.0834

$ws_ARM.Range("H77").Value = 1293.762
$ws_ARM.Range("I77").Value = 1098.6666
$ws_ARM.Range("J77").Value = 1440.0834
$ws_ARM.Range("K77").Value = 5493.333000000001
$ws_ARM.Range("L77").Value = 7200.416999999999
$ws_ARM.Range("M77").Value = -1125.333000000001
$ws_ARM.Range("N77").Value = -15936.417

$ws_ARM.Range("H92").Value = 74698
$ws_ARM.Range("I92").Value = 0
$ws_ARM.Range("J92").Value = 74698
$ws_ARM.Range("K92").Value = 0
$ws_ARM.Range("L92").Value = 74698
$ws_ARM.Range("N92").Value = -79690

$ws_ARM.Range("H132").Value = 3027.3572
$ws_ARM.Range("I132").Value = 1798.125
$ws_ARM.Range("J132").Value = 4666.3335
$ws_ARM.Range("K132").Value = 5394.375
$ws_ARM.Range("L132").Value = 13999.0005
$ws_ARM.Range("M132").Value = -2864.375
$ws_ARM.Range("N132").Value = -19059.0005

$ws_ARM.Range("H134").Value = 52000
$ws_ARM.Range("I134").Value = 0
$ws_ARM.Range("J134").Value = 52000
$ws_ARM.Range("K134").Value = 0
$ws_ARM.Range("L134").Value = 52000
$ws_ARM.Range("N134").Value = -62140

$ws_ARM.Range("H139").Value = 85000
$ws_ARM.Range("I139").Value = 0
$ws_ARM.Range("J139").Value = 85000
$ws_ARM.Range("K139").Value = 0
$ws_ARM.Range("L139").Value = 85000
$ws_ARM.Range("N139").Value = -95280

$ws_ARM.Range("H141").Value = 35913.223
$ws_ARM.Range("I141").Value = 0
$ws_ARM.Range("J141").Value = 35913.223
$ws_ARM.Range("K141").Value = 0
$ws_ARM.Range("L141").Value = 35913.223
$ws_ARM.Range("N141").Value = -46273.223

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 2351.772
$ws_CRP.Range("I31").Value = 1859.4651
$ws_CRP.Range("J31").Value = 3863.8572
$ws_CRP.Range("K31").Value = 1859.4651
$ws_CRP.Range("L31").Value = 3863.8572
$ws_CRP.Range("M31").Value = -1564.4651
$ws_CRP.Range("N31").Value = -4453.8572

$ws_CRP.Range("H33").Value = 0
$ws_CRP.Range("I33").Value = 0
$ws_CRP.Range("J33").Value = 0
$ws_CRP.Range("K33").Value = 0
$ws_CRP.Range("L33").ClearContents()
$ws_CRP.Range("N33").Value = 0

$ws_CRP.Range("H34").Value = 2351.772
$ws_CRP.Range("I34").Value = 1859.4651
$ws_CRP.Range("J34").Value = 3863.8572
$ws_CRP.Range("K34").Value = 1859.4651
$ws_CRP.Range("L34").Value = 3863.8572
$ws_CRP.Range("M34").Value = -1657.4651
$ws_CRP.Range("N34").Value = -4267.8572

$ws_CRP.Range("H37").Value = 15000
$ws_CRP.Range("I37").Value = 0
$ws_CRP.Range("J37").Value = 15000
$ws_CRP.Range("K37").Value = 0
$ws_CRP.Range("L37").Value = 15000
$ws_CRP.Range("N37").Value = -15214

$ws_CRP.Range("H58").Value = 1422.9524
$ws_CRP.Range("I58").Value = 1330.6316
$ws_CRP.Range("J58").Value = 2300
$ws_CRP.Range("K58").Value = 1330.6316
$ws_CRP.Range("L58").Value = 2300
$ws_CRP.Range("M58").Value = -1127.6316
$ws_CRP.Range("N58").Value = -2706

$ws_CRP.Range("H132").Value = 1984.3125
$ws_CRP.Range("I132").Value = 1468.3636
$ws_CRP.Range("J132").Value = 3119.4
$ws_CRP.Range("K132").Value = 4405.0908
$ws_CRP.Range("L132").Value = 9358.200000000001
$ws_CRP.Range("M132").Value = -1875.0908
$ws_CRP.Range("N132").Value = -14418.2

$ws_CRP.Range("H133").Value = 0
$ws_CRP.Range("I133").Value = 0
$ws_CRP.Range("J133").Value = 0
$ws_CRP.Range("K133").Value = 0
$ws_CRP.Range("L133").ClearContents()
$ws_CRP.Range("N133").Value = 0

$ws_CRP.Range("H136").Value = 1422.9524
$ws_CRP.Range("I136").Value = 1330.6316
$ws_CRP.Range("J136").Value = 2300
$ws_CRP.Range("K136").Value = 3991.8948
$ws_CRP.Range("L136").Value = 6900
$ws_CRP.Range("M136").Value = -1441.8948
$ws_CRP.Range("N136").Value = -12000

$ws_CRP.Range("H140").Value = 0
$ws_CRP.Range("I140").Value = 0
$ws_CRP.Range("J140").Value = 0
$ws_CRP.Range("K140").Value = 0
$ws_CRP.Range("L140").ClearContents()
$ws_CRP.Range("N140").Value = 0

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 1508.3846
$ws_CUL.Range("I5").Value = 1918.6
$ws_CUL.Range("J5").Value = 949
$ws_CUL.Range("K5").Value = 5755.799999999999
$ws_CUL.Range("L5").Value = 2847
$ws_CUL.Range("M5").Value = -5643.799999999999
$ws_CUL.Range("N5").Value = -3071

$ws_CUL.Range("H34").Value = 906.38464
$ws_CUL.Range("I34").Value = 290
$ws_CUL.Range("J34").Value = 957.75
$ws_CUL.Range("K34").Value = 870
$ws_CUL.Range("L34").Value = 2873.25
$ws_CUL.Range("M34").Value = -786
$ws_CUL.Range("N34").Value = -3041.25

$ws_CUL.Range("H113").Value = 476830.84
$ws_CUL.Range("I113").Value = 667285.25
$ws_CUL.Range("J113").Value = 694.8333
$ws_CUL.Range("K113").Value = 2001855.75
$ws_CUL.Range("L113").Value = 2084.4999
$ws_CUL.Range("M113").Value = -1999685.75
$ws_CUL.Range("N113").Value = -6424.4999

$ws_CUL.Range("H131").Value = 12988389
$ws_CUL.Range("I131").Value = 522
$ws_CUL.Range("J131").Value = 13890324
$ws_CUL.Range("K131").Value = 1566
$ws_CUL.Range("L131").Value = 41670972
$ws_CUL.Range("M131").Value = 3474
$ws_CUL.Range("N131").Value = -41681052

$ws_CUL.Range("H135").Value = 1508.3846
$ws_CUL.Range("I135").Value = 1918.6
$ws_CUL.Range("J135").Value = 949
$ws_CUL.Range("K135").Value = 17267.4
$ws_CUL.Range("L135").Value = 8541
$ws_CUL.Range("M135").Value = -14732.4
$ws_CUL.Range("N135").Value = -13611

$ws_CUL.Range("H137").Value = 2948.75
$ws_CUL.Range("I137").Value = 1514.5
$ws_CUL.Range("J137").Value = 4383
$ws_CUL.Range("K137").Value = 4543.5
$ws_CUL.Range("L137").Value = 13149
$ws_CUL.Range("M137").Value = 556.5
$ws_CUL.Range("N137").Value = -23349

$ws_CUL.Range("H139").Value = 1681.7693
$ws_CUL.Range("I139").Value = 1103.3334
$ws_CUL.Range("J139").Value = 2983.25
$ws_CUL.Range("K139").Value = 3310.0002
$ws_CUL.Range("L139").Value = 8949.75
$ws_CUL.Range("M139").Value = 1829.9998
$ws_CUL.Range("N139").Value = -19229.75

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H34").Value = 36633.332
$ws_LTW.Range("I34").Value = 44950
$ws_LTW.Range("J34").Value = 20000
$ws_LTW.Range("K34").Value = 44950
$ws_LTW.Range("L34").Value = 20000
$ws_LTW.Range("M34").Value = -44778
$ws_LTW.Range("N34").Value = -20344

$ws_LTW.Range("H132").Value = 5095.921
$ws_LTW.Range("I132").Value = 4508.375
$ws_LTW.Range("J132").Value = 6103.143
$ws_LTW.Range("K132").Value = 13525.125
$ws_LTW.Range("L132").Value = 18309.429
$ws_LTW.Range("M132").Value = -10995.125
$ws_LTW.Range("N132").Value = -23369.429

$ws_LTW.Range("H133").Value = 60328
$ws_LTW.Range("I133").Value = 0
$ws_LTW.Range("J133").Value = 60328
$ws_LTW.Range("K133").Value = 0
$ws_LTW.Range("L133").Value = 60328
$ws_LTW.Range("N133").Value = -65388

$ws_LTW.Range("H137").Value = 74852.375
$ws_LTW.Range("I137").Value = 49000
$ws_LTW.Range("J137").Value = 83469.836
$ws_LTW.Range("K137").Value = 49000
$ws_LTW.Range("L137").Value = 83469.836
$ws_LTW.Range("M137").Value = -43900
$ws_LTW.Range("N137").Value = -93669.836

$ws_LTW.Range("H138").Value = 50000
$ws_LTW.Range("I138").Value = 0
$ws_LTW.Range("J138").Value = 50000
$ws_LTW.Range("K138").Value = 0
$ws_LTW.Range("L138").Value = 50000
$ws_LTW.Range("N138").Value = -60280

$ws_LTW.Range("H141").Value = 45000
$ws_LTW.Range("I141").Value = 0
$ws_LTW.Range("J141").Value = 45000
$ws_LTW.Range("K141").Value = 0
$ws_LTW.Range("L141").Value = 45000
$ws_LTW.Range("N141").Value = -55360

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H100").Value = 10692.667
$ws_WVR.Range("I100").Value = 17561.666
$ws_WVR.Range("J100").Value = 1534
$ws_WVR.Range("K100").Value = 35123.332
$ws_WVR.Range("L100").Value = 3068
$ws_WVR.Range("M100").Value = -34582.332
$ws_WVR.Range("N100").Value = -4150

$ws_WVR.Range("H132").Value = 2215.7727
$ws_WVR.Range("I132").Value = 1626
$ws_WVR.Range("J132").Value = 2805.5454
$ws_WVR.Range("K132").Value = 4878
$ws_WVR.Range("L132").Value = 8416.636200000001
$ws_WVR.Range("M132").Value = -2348
$ws_WVR.Range("N132").Value = -13476.6362

$ws_WVR.Range("H135").Value = 35500
$ws_WVR.Range("I135").Value = 0
$ws_WVR.Range("J135").Value = 35500
$ws_WVR.Range("K135").Value = 0
$ws_WVR.Range("L135").Value = 35500
$ws_WVR.Range("N135").Value = -45640

$ws_WVR.Range("H136").Value = 1322.9642
$ws_WVR.Range("I136").Value = 1356
$ws_WVR.Range("J136").Value = 1124.75
$ws_WVR.Range("K136").Value = 4068
$ws_WVR.Range("L136").Value = 3374.25
$ws_WVR.Range("M136").Value = -1518
$ws_WVR.Range("N136").Value = -8474.25

$ws_WVR.Range("H137").Value = 49827.5
$ws_WVR.Range("I137").Value = 0
$ws_WVR.Range("J137").Value = 49827.5
$ws_WVR.Range("K137").Value = 0
$ws_WVR.Range("L137").Value = 49827.5
$ws_WVR.Range("N137").Value = -60027.5

$ws_WVR.Range("H138").Value = 62981.332
$ws_WVR.Range("I138").Value = 0
$ws_WVR.Range("J138").Value = 62981.332
$ws_WVR.Range("K138").Value = 0
$ws_WVR.Range("L138").Value = 62981.332
$ws_WVR.Range("N138").Value = -73261.33199999999

$ws_WVR.Range("H139").Value = 50996
$ws_WVR.Range("I139").Value = 0
$ws_WVR.Range("J139").Value = 50996
$ws_WVR.Range("K139").Value = 0
$ws_WVR.Range("L139").Value = 50996
$ws_WVR.Range("N139").Value = -61276

$ws_WVR.Range("H140").Value = 49333.332
$ws_WVR.Range("I140").Value = 0
$ws_WVR.Range("J140").Value = 49333.332
$ws_WVR.Range("K140").Value = 0
$ws_WVR.Range("L140").Value = 49333.332
$ws_WVR.Range("N140").Value = -59693.332
